$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text
$ws.Range("A1").Value = "Datos actualizados a 2 de Agosto de 2020 a las 09:30"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 4764522
$ws.Range("C4").Value = 204
$ws.Range("D4").Value = 2363001
$ws.Range("E4").Value = 2243620
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 157901

# Row 37: Ucrania
$ws.Range("A37").Value = "Ucrania"
$ws.Range("B37").Value = 72168
$ws.Range("C37").Value = 1112
$ws.Range("D37").Value = 39543
$ws.Range("E37").Value = 30900
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 16
$ws.Range("H37").Value = 1725

# Row 38: Republica Dominicana
$ws.Range("A38").Value = "Republica Dominicana"
$ws.Range("B38").Value = 71415
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 37509
$ws.Range("E38").Value = 32736
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 1170

# Row 53: Armenia
$ws.Range("A53").Value = "Armenia"
$ws.Range("B53").Value = 39050
$ws.Range("C53").Value = 209
$ws.Range("D53").Value = 29750
$ws.Range("E53").Value = 8546
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 5
$ws.Range("H53").Value = 754

# Row 106: Hungria
$ws.Range("A106").Value = "Hungria"
$ws.Range("B106").Value = 4535
$ws.Range("C106").Value = 9
$ws.Range("D106").Value = 3389
$ws.Range("E106").Value = 549
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 597

# Row 128: Estonia
$ws.Range("A128").Value = "Estonia"
$ws.Range("B128").Value = 2079
$ws.Range("C128").Value = 7
$ws.Range("D128").Value = 1934
$ws.Range("E128").Value = 76
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 69

# Row 140: Letonia
$ws.Range("A140").Value = "Letonia"
$ws.Range("B140").Value = 1243
$ws.Range("C140").Value = 5
$ws.Range("D140").Value = 1052
$ws.Range("E140").Value = 159
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 32

# Row 143: Georgia
$ws.Range("A143").Value = "Georgia"
$ws.Range("B143").Value = 1177
$ws.Range("C143").Value = 6
$ws.Range("D143").Value = 955
$ws.Range("E143").Value = 205
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 17

# Row 144: Uganda
$ws.Range("A144").Value = "Uganda"
$ws.Range("B144").Value = 1176
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 1045
$ws.Range("E144").Value = 127
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 4

# Row 162: Vietnam
$ws.Range("A162").Value = "Vietnam"
$ws.Range("B162").Value = 590
$ws.Range("C162").Value = 0
$ws.Range("D162").Value = 373
$ws.Range("E162").Value = 212
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 2
$ws.Range("H162").Value = 5

# Row 165: Taiwan
$ws.Range("A165").Value = "Taiwan"
$ws.Range("B165").Value = 475
$ws.Range("C165").Value = 1
$ws.Range("D165").Value = 441
$ws.Range("E165").Value = 27
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 7

# Row 211: Bonaire, San Eustaquio y Saba
$ws.Range("A211").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("B211").Value = 13
$ws.Range("C211").Value = 2
$ws.Range("D211").Value = 7
$ws.Range("E211").Value = 6
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0

# Row 212: Islas Malvinas
$ws.Range("A212").Value = "Islas Malvinas"
$ws.Range("B212").Value = 13
$ws.Range("C212").Value = 0
$ws.Range("D212").Value = 13
$ws.Range("E212").Value = 0
$ws.Range("F212").Value = 0
$ws.Range("G212").Value = 0
$ws.Range("H212").Value = 0

# Row 213: Montserrat
$ws.Range("A213").Value = "Montserrat"
$ws.Range("B213").Value = 12
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 10
$ws.Range("E213").Value = 1
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

# Row 214: Santa Sede
$ws.Range("A214").Value = "Santa Sede"
$ws.Range("B214").Value = 12
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0
